$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 536, pushing the existing
# rows 536:542 down to 539:545 (same as the source diff shows).
$ws.Rows("536:538").Insert()

# New row 536 - Morada(o), 1a (cosecha), malla 18 kilos, O'Higgins
$ws.Range("A536").Value = 7
$ws.Range("B536").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C536").Value = "Ñuble"
$ws.Range("D536").Value = 44656
$ws.Range("E536").Value = 16
$ws.Range("F536").Value = 100112004
$ws.Range("G536").Value = "Cebolla"
$ws.Range("H536").Value = "Morada(o)"
$ws.Range("I536").Value = "1a (cosecha)"
$ws.Range("J536").Value = 120
$ws.Range("K536").Value = 8000
$ws.Range("L536").Value = 8500
$ws.Range("M536").Value = 8250
$ws.Range("N536").Value = "`$/malla 18 kilos"
$ws.Range("O536").Value = "Región de O'Higgins"
$ws.Range("P536").Value = 458
$ws.Range("Q536").Value = 18
$ws.Range("R536").Value = "Hortaliza"

# New row 537 - Sin especificar, 1a (cosecha), malla 25 kilos, O'Higgins
$ws.Range("A537").Value = 7
$ws.Range("B537").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C537").Value = "Ñuble"
$ws.Range("D537").Value = 44656
$ws.Range("E537").Value = 16
$ws.Range("F537").Value = 100112004
$ws.Range("G537").Value = "Cebolla"
$ws.Range("H537").Value = "Sin especificar"
$ws.Range("I537").Value = "1a (cosecha)"
$ws.Range("J537").Value = 200
$ws.Range("K537").Value = 5000
$ws.Range("L537").Value = 5500
$ws.Range("M537").Value = 5250
$ws.Range("N537").Value = "`$/malla 25 kilos"
$ws.Range("O537").Value = "Región de O'Higgins"
$ws.Range("P537").Value = 210
$ws.Range("Q537").Value = 25
$ws.Range("R537").Value = "Hortaliza"

# New row 538 - Sin especificar, 2a (cosecha), malla 25 kilos, Maule
$ws.Range("A538").Value = 7
$ws.Range("B538").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C538").Value = "Ñuble"
$ws.Range("D538").Value = 44656
$ws.Range("E538").Value = 16
$ws.Range("F538").Value = 100112004
$ws.Range("G538").Value = "Cebolla"
$ws.Range("H538").Value = "Sin especificar"
$ws.Range("I538").Value = "2a (cosecha)"
$ws.Range("J538").Value = 120
$ws.Range("K538").Value = 4000
$ws.Range("L538").Value = 4500
$ws.Range("M538").Value = 4250
$ws.Range("N538").Value = "`$/malla 25 kilos"
$ws.Range("O538").Value = "Región del Maule"
$ws.Range("P538").Value = 170
$ws.Range("Q538").Value = 25
$ws.Range("R538").Value = "Hortaliza"
